$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.541.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.674.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.84%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5285'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.32%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2680'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06390'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07780'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.493'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.667.71'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5575'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅8336'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.526.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.09%  '

$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.769'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.313'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1274'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.32%  '

$ws.Range("E26").Value = '  +0.31%  '

$ws.Range("E27").Value = '  +2.99%  '

$ws.Range("E28").Value = '  +2.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06215'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.292'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.632'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.436'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.26%  '

$ws.Range("E33").Value = '  +2.38%  '

$ws.Range("E34").Value = '  +1.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6100'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.12%  '

$ws.Range("E36").Value = '  +0.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.780'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.95%  '

$ws.Range("E38").Value = '  +0.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.063'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.088.66'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8608'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.06%  '

$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.818.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.41%  '

$ws.Range("E45").Value = '  +2.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.007'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.527'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.109'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05196'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.020'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.03%  '

